$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $orig = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = $orig
}

Set-TextValue $ws.Cells.Item(2, 4) "310.60"
Set-TextValue $ws.Cells.Item(2, 5) "2.06%"
Set-TextValue $ws.Cells.Item(3, 4) "38.87"
Set-TextValue $ws.Cells.Item(3, 5) "8.63%"
Set-TextValue $ws.Cells.Item(4, 4) "5.116"
Set-TextValue $ws.Cells.Item(4, 5) "1.59%"
Set-TextValue $ws.Cells.Item(5, 4) "0.08183"
Set-TextValue $ws.Cells.Item(5, 5) "2.81%"
Set-TextValue $ws.Cells.Item(6, 4) "2.014"
Set-TextValue $ws.Cells.Item(6, 5) "7.78%"
Set-TextValue $ws.Cells.Item(7, 2) "KuCoinToken"
Set-TextValue $ws.Cells.Item(7, 3) "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Cells.Item(7, 4) "7.927"
Set-TextValue $ws.Cells.Item(7, 5) "2.03%"
Set-TextValue $ws.Cells.Item(8, 2) "MXToken"
Set-TextValue $ws.Cells.Item(8, 3) "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(8, 4) "0.9348"
Set-TextValue $ws.Cells.Item(8, 5) "1.56%"
Set-TextValue $ws.Cells.Item(9, 2) "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Cells.Item(9, 3) "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Cells.Item(9, 4) "0.1409"
Set-TextValue $ws.Cells.Item(9, 5) "5.12%"
Set-TextValue $ws.Cells.Item(10, 2) "WazirX"
Set-TextValue $ws.Cells.Item(10, 3) "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Cells.Item(10, 4) "0.1962"
Set-TextValue $ws.Cells.Item(10, 5) "3.78%"
Set-TextValue $ws.Cells.Item(11, 2) "MandalaExchangeToken"
Set-TextValue $ws.Cells.Item(11, 3) "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Cells.Item(11, 4) "0.09185"
Set-TextValue $ws.Cells.Item(11, 5) "0.43%"
Set-TextValue $ws.Cells.Item(12, 2) "BitrueCoin"
Set-TextValue $ws.Cells.Item(12, 3) "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Cells.Item(12, 4) "0.03458"
Set-TextValue $ws.Cells.Item(12, 5) "0.67%"
Set-TextValue $ws.Cells.Item(13, 2) "BitMartToken"
Set-TextValue $ws.Cells.Item(13, 3) "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Cells.Item(13, 4) "0.09841"
Set-TextValue $ws.Cells.Item(13, 5) "-0.07%"
Set-TextValue $ws.Cells.Item(14, 2) "BitForexToken"
Set-TextValue $ws.Cells.Item(14, 3) "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Cells.Item(14, 4) "0.001411"
Set-TextValue $ws.Cells.Item(14, 5) "0.38%"
Set-TextValue $ws.Cells.Item(15, 2) "TigerCash"
Set-TextValue $ws.Cells.Item(15, 3) "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Cells.Item(15, 4) "0.005866"
Set-TextValue $ws.Cells.Item(15, 5) "-3.38%"
Set-TextValue $ws.Cells.Item(16, 2) "LEO"
Set-TextValue $ws.Cells.Item(16, 3) "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Cells.Item(16, 4) "3.570"
Set-TextValue $ws.Cells.Item(16, 5) "-4.26%"
Set-TextValue $ws.Cells.Item(17, 2) "GateToken"
Set-TextValue $ws.Cells.Item(17, 3) "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Cells.Item(17, 4) "4.197"
Set-TextValue $ws.Cells.Item(17, 5) "1.97%"
Set-TextValue $ws.Cells.Item(19, 5) "0.20%"
Set-TextValue $ws.Cells.Item(20, 5) "-0.16%"
Set-TextValue $ws.Cells.Item(21, 4) "4.829"
Set-TextValue $ws.Cells.Item(21, 5) "-6.46%"
Set-TextValue $ws.Cells.Item(22, 5) "5.16%"
Set-TextValue $ws.Cells.Item(23, 4) "0.04475"
Set-TextValue $ws.Cells.Item(23, 5) "1.34%"
Set-TextValue $ws.Cells.Item(24, 4) "0.001240"
Set-TextValue $ws.Cells.Item(24, 5) "1.10%"
Set-TextValue $ws.Cells.Item(25, 5) "-9.72%"
Set-TextValue $ws.Cells.Item(27, 4) "0.0001302"
Set-TextValue $ws.Cells.Item(27, 5) "0.28%"
Set-TextValue $ws.Cells.Item(39, 4) "0.02117"
Set-TextValue $ws.Cells.Item(39, 5) "8.71%"
Set-TextValue $ws.Cells.Item(40, 4) "0.05182"
Set-TextValue $ws.Cells.Item(40, 5) "-5.26%"
Set-TextValue $ws.Cells.Item(41, 4) "0.007461"
Set-TextValue $ws.Cells.Item(41, 5) "-1.83%"
Set-TextValue $ws.Cells.Item(42, 4) "0.009964"
Set-TextValue $ws.Cells.Item(42, 5) "-1.62%"
Set-TextValue $ws.Cells.Item(43, 4) "0.1368"
Set-TextValue $ws.Cells.Item(43, 5) "1.19%"
Set-TextValue $ws.Cells.Item(44, 4) "0.002134"
Set-TextValue $ws.Cells.Item(44, 5) "-0.65%"
Set-TextValue $ws.Cells.Item(45, 4) "0.009759"
Set-TextValue $ws.Cells.Item(45, 5) "-4.06%"
Set-TextValue $ws.Cells.Item(46, 4) "0.00006341"
Set-TextValue $ws.Cells.Item(46, 5) "3.82%"
Set-TextValue $ws.Cells.Item(47, 5) "0.32%"
Set-TextValue $ws.Cells.Item(48, 5) "-0.24%"
Set-TextValue $ws.Cells.Item(49, 4) "0.001603"
Set-TextValue $ws.Cells.Item(49, 5) "-3.23%"
Set-TextValue $ws.Cells.Item(50, 5) "0.32%"
Set-TextValue $ws.Cells.Item(51, 4) "0.0002004"
Set-TextValue $ws.Cells.Item(51, 5) "0.32%"
